# Append the latest API call results (two new POST/GET request pairs)
# to the "Results" sheet, and widen the ID column to fit them.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Operation, ID, Name, Year, Price, CPU, HardDisk, CreatedAt, Timestamp
$rows = @(
    @("POST", "ff8081819782e69e019908606a771630", "Apple MacBook Pro 20 Max", 2021, 189.99, "Intel Core i8", "8TB", "2025-09-02T03:02:40.503+00:00", "2025-09-02 08:32:40"),
    @("GET",  "ff8081819782e69e019908606a771630", "Apple MacBook Pro 20 Max", 2021, 189.99, "Intel Core i8", "8TB", $null, "2025-09-02 08:32:40"),
    @("POST", "ff8081819782e69e0199086097711631", "Apple MacBook Pro 20 Max", 2021, 189.99, "Intel Core i8", "8TB", "2025-09-02T03:02:52.017+00:00", "2025-09-02 08:32:51"),
    @("GET",  "ff8081819782e69e0199086097711631", "Apple MacBook Pro 20 Max", 2021, 189.99, "Intel Core i8", "8TB", $null, "2025-09-02 08:32:52")
)

$r = 4
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    # GET rows have no CreatedAt (request never hit that stage)
    if ($row[7] -ne $null) {
        $ws.Cells.Item($r, 8).Value = $row[7]
    }
    $ws.Cells.Item($r, 9).Value = $row[8]
    $r = $r + 1
}

# Widen the ID column now that it holds the new, longer record ids
$ws.Columns.Item(2).ColumnWidth = 52
